$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
  # row 18: H18=2484, I18=2484, K18=2484, M18=-2200
  $ws.Range("H18").Value = 2484
  $ws.Range("I18").Value = 2484
  $ws.Range("K18").Value = 2484
  $ws.Range("M18").Value = -2200
  # row 43: H43=9351.454, I43=14999, J43=8096.4443, K43=14999, L43=8096.4443, M43=-14930, N43=-8234.444299999999
  $ws.Range("H43").Value = 9351.454
  $ws.Range("I43").Value = 14999
  $ws.Range("J43").Value = 8096.4443
  $ws.Range("K43").Value = 14999
  $ws.Range("L43").Value = 8096.4443
  $ws.Range("M43").Value = -14930
  $ws.Range("N43").Value = -8234.444299999999
  # row 86: H86=6093.9, I86=4179.6665, J86=6914.2856, K86=4179.6665, L86=6914.2856, M86=-3056.6665, N86=-9160.285599999999
  $ws.Range("H86").Value = 6093.9
  $ws.Range("I86").Value = 4179.6665
  $ws.Range("J86").Value = 6914.2856
  $ws.Range("K86").Value = 4179.6665
  $ws.Range("L86").Value = 6914.2856
  $ws.Range("M86").Value = -3056.6665
  $ws.Range("N86").Value = -9160.285599999999
  # row 89: H89=6093.9, I89=4179.6665, J89=6914.2856, K89=20898.3325, L89=34571.428, M89=-15282.3325, N89=-45803.428
  $ws.Range("H89").Value = 6093.9
  $ws.Range("I89").Value = 4179.6665
  $ws.Range("J89").Value = 6914.2856
  $ws.Range("K89").Value = 20898.3325
  $ws.Range("L89").Value = 34571.428
  $ws.Range("M89").Value = -15282.3325
  $ws.Range("N89").Value = -45803.428
  # row 96: H96=1651, J96=2364.75, L96=7094.25, N96=-9840.25
  $ws.Range("H96").Value = 1651
  $ws.Range("J96").Value = 2364.75
  $ws.Range("L96").Value = 7094.25
  $ws.Range("N96").Value = -9840.25
  # row 129: H129=1482.5, I129=710.1667, J129=2061.75, K129=2130.5001, L129=6185.25, M129=2869.4999, N129=-16185.25
  $ws.Range("H129").Value = 1482.5
  $ws.Range("I129").Value = 710.1667
  $ws.Range("J129").Value = 2061.75
  $ws.Range("K129").Value = 2130.5001
  $ws.Range("L129").Value = 6185.25
  $ws.Range("M129").Value = 2869.4999
  $ws.Range("N129").Value = -16185.25
  # row 137: H137=1572.3658, I137=1065.8889, K137=3197.6667, M137=-647.6666999999998
  $ws.Range("H137").Value = 1572.3658
  $ws.Range("I137").Value = 1065.8889
  $ws.Range("K137").Value = 3197.6667
  $ws.Range("M137").Value = -647.6666999999998
  # row 138: H138=4855.904, I138=1341.8695, J138=7642.8965, K138=4025.6085, L138=22928.6895, M138=1114.3915, N138=-33208.6895
  $ws.Range("H138").Value = 4855.904
  $ws.Range("I138").Value = 1341.8695
  $ws.Range("J138").Value = 7642.8965
  $ws.Range("K138").Value = 4025.6085
  $ws.Range("L138").Value = 22928.6895
  $ws.Range("M138").Value = 1114.3915
  $ws.Range("N138").Value = -33208.6895
  # row 141: H141=2754, I141=2224.4443, K141=6673.3329, M141=-1493.3329
  $ws.Range("H141").Value = 2754
  $ws.Range("I141").Value = 2224.4443
  $ws.Range("K141").Value = 6673.3329
  $ws.Range("M141").Value = -1493.3329

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
  # row 2: H2=1512.0303, I2=1246.1, J2=1921.1538, K2=1246.1, L2=1921.1538, M2=-1133.1, N2=-2147.1538
  $ws.Range("H2").Value = 1512.0303
  $ws.Range("I2").Value = 1246.1
  $ws.Range("J2").Value = 1921.1538
  $ws.Range("K2").Value = 1246.1
  $ws.Range("L2").Value = 1921.1538
  $ws.Range("M2").Value = -1133.1
  $ws.Range("N2").Value = -2147.1538
  # row 32: H32=2987.1765, I32=3427.6365, K32=3427.6365, M32=-3140.6365
  $ws.Range("H32").Value = 2987.1765
  $ws.Range("I32").Value = 3427.6365
  $ws.Range("K32").Value = 3427.6365
  $ws.Range("M32").Value = -3140.6365
  # row 61: H61=5765.625, I61=2557.8276, J61=9211.037, K61=2557.8276, L61=9211.037, M61=-2345.8276, N61=-9635.037
  $ws.Range("H61").Value = 5765.625
  $ws.Range("I61").Value = 2557.8276
  $ws.Range("J61").Value = 9211.037
  $ws.Range("K61").Value = 2557.8276
  $ws.Range("L61").Value = 9211.037
  $ws.Range("M61").Value = -2345.8276
  $ws.Range("N61").Value = -9635.037
  # row 63: H63=3936.5833, I63=3936.5833, K63=3936.5833, M63=-3250.5833
  $ws.Range("H63").Value = 3936.5833
  $ws.Range("I63").Value = 3936.5833
  $ws.Range("K63").Value = 3936.5833
  $ws.Range("M63").Value = -3250.5833
  # row 66: H66=3936.5833, I66=3936.5833, K66=19682.9165, M66=-16250.9165
  $ws.Range("H66").Value = 3936.5833
  $ws.Range("I66").Value = 3936.5833
  $ws.Range("K66").Value = 19682.9165
  $ws.Range("M66").Value = -16250.9165
  # row 116: H116=1512.0303, I116=1246.1, J116=1921.1538, K116=1246.1, L116=1921.1538, M116=1047.9, N116=-6509.1538
  $ws.Range("H116").Value = 1512.0303
  $ws.Range("I116").Value = 1246.1
  $ws.Range("J116").Value = 1921.1538
  $ws.Range("K116").Value = 1246.1
  $ws.Range("L116").Value = 1921.1538
  $ws.Range("M116").Value = 1047.9
  $ws.Range("N116").Value = -6509.1538
  # row 136: H136=5765.625, I136=2557.8276, J136=9211.037, K136=7673.4828, L136=27633.111, M136=-5123.4828, N136=-32733.111
  $ws.Range("H136").Value = 5765.625
  $ws.Range("I136").Value = 2557.8276
  $ws.Range("J136").Value = 9211.037
  $ws.Range("K136").Value = 7673.4828
  $ws.Range("L136").Value = 27633.111
  $ws.Range("M136").Value = -5123.4828
  $ws.Range("N136").Value = -32733.111

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
  # row 3: H3=1512.0303, I3=1246.1, J3=1921.1538, K3=1246.1, L3=1921.1538, M3=-1132.1, N3=-2149.1538
  $ws.Range("H3").Value = 1512.0303
  $ws.Range("I3").Value = 1246.1
  $ws.Range("J3").Value = 1921.1538
  $ws.Range("K3").Value = 1246.1
  $ws.Range("L3").Value = 1921.1538
  $ws.Range("M3").Value = -1132.1
  $ws.Range("N3").Value = -2149.1538
  # row 134: H134=1718.5312, I134=1525.6552, K134=4576.9656, M134=-2041.9656
  $ws.Range("H134").Value = 1718.5312
  $ws.Range("I134").Value = 1525.6552
  $ws.Range("K134").Value = 4576.9656
  $ws.Range("M134").Value = -2041.9656

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
  # row 6: H6=6575885, I6=7671699.5, K6=7671699.5, M6=-7671586.5
  $ws.Range("H6").Value = 6575885
  $ws.Range("I6").Value = 7671699.5
  $ws.Range("K6").Value = 7671699.5
  $ws.Range("M6").Value = -7671586.5
  # row 22: H22=499.625, I22=499.83334, J22=499, K22=499.83334, L22=499, M22=-149.83334, N22=-1199
  $ws.Range("H22").Value = 499.625
  $ws.Range("I22").Value = 499.83334
  $ws.Range("J22").Value = 499
  $ws.Range("K22").Value = 499.83334
  $ws.Range("L22").Value = 499
  $ws.Range("M22").Value = -149.83334
  $ws.Range("N22").Value = -1199
  # row 31: H31=3062.8, I31=1983.9524, K31=1983.9524, M31=-1688.9524
  $ws.Range("H31").Value = 3062.8
  $ws.Range("I31").Value = 1983.9524
  $ws.Range("K31").Value = 1983.9524
  $ws.Range("M31").Value = -1688.9524
  # row 34: H34=3062.8, I34=1983.9524, K34=1983.9524, M34=-1781.9524
  $ws.Range("H34").Value = 3062.8
  $ws.Range("I34").Value = 1983.9524
  $ws.Range("K34").Value = 1983.9524
  $ws.Range("M34").Value = -1781.9524
  # row 58: H58=2523.8333, I58=2440, J58=2663.5557, K58=2440, L58=2663.5557, M58=-2237, N58=-3069.5557
  $ws.Range("H58").Value = 2523.8333
  $ws.Range("I58").Value = 2440
  $ws.Range("J58").Value = 2663.5557
  $ws.Range("K58").Value = 2440
  $ws.Range("L58").Value = 2663.5557
  $ws.Range("M58").Value = -2237
  $ws.Range("N58").Value = -3069.5557
  # row 136: H136=2523.8333, I136=2440, J136=2663.5557, K136=7320, L136=7990.6671, M136=-4770, N136=-13090.6671
  $ws.Range("H136").Value = 2523.8333
  $ws.Range("I136").Value = 2440
  $ws.Range("J136").Value = 2663.5557
  $ws.Range("K136").Value = 7320
  $ws.Range("L136").Value = 7990.6671
  $ws.Range("M136").Value = -4770
  $ws.Range("N136").Value = -13090.6671

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
  # row 43: H43=5959.4, J43=7450, L43=22350, N43=-22578
  $ws.Range("H43").Value = 5959.4
  $ws.Range("J43").Value = 7450
  $ws.Range("L43").Value = 22350
  $ws.Range("N43").Value = -22578
  # row 114: H114=1669.9166, I114=1371.75, J114=1819, K114=4115.25, L114=5457, M114=-861.25, N114=-11965
  $ws.Range("H114").Value = 1669.9166
  $ws.Range("I114").Value = 1371.75
  $ws.Range("J114").Value = 1819
  $ws.Range("K114").Value = 4115.25
  $ws.Range("L114").Value = 5457
  $ws.Range("M114").Value = -861.25
  $ws.Range("N114").Value = -11965
  # row 132: H132=1366.1666, I132=1099.25, J132=1900, K132=9893.25, L132=17100, M132=-7363.25, N132=-22160
  $ws.Range("H132").Value = 1366.1666
  $ws.Range("I132").Value = 1099.25
  $ws.Range("J132").Value = 1900
  $ws.Range("K132").Value = 9893.25
  $ws.Range("L132").Value = 17100
  $ws.Range("M132").Value = -7363.25
  $ws.Range("N132").Value = -22160

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
  # row 70: H70=29400, I70=8800, K70=8800, M70=-8530
  $ws.Range("H70").Value = 29400
  $ws.Range("I70").Value = 8800
  $ws.Range("K70").Value = 8800
  $ws.Range("M70").Value = -8530
  # row 73: H73=29400, I73=8800, K73=8800, M73=-7864
  $ws.Range("H73").Value = 29400
  $ws.Range("I73").Value = 8800
  $ws.Range("K73").Value = 8800
  $ws.Range("M73").Value = -7864
  # row 102: H102=21143.54, I102=2005.5834, K102=2005.5834, M102=-383.5834
  $ws.Range("H102").Value = 21143.54
  $ws.Range("I102").Value = 2005.5834
  $ws.Range("K102").Value = 2005.5834
  $ws.Range("M102").Value = -383.5834
  # row 113: H113=2915.1365, I113=3034.1765, J113=2510.4, K113=3034.1765, L113=2510.4, M113=-864.1765, N113=-6850.4
  $ws.Range("H113").Value = 2915.1365
  $ws.Range("I113").Value = 3034.1765
  $ws.Range("J113").Value = 2510.4
  $ws.Range("K113").Value = 3034.1765
  $ws.Range("L113").Value = 2510.4
  $ws.Range("M113").Value = -864.1765
  $ws.Range("N113").Value = -6850.4
  # row 132: H132=3644.2917, J132=1990.1428, L132=5970.428400000001, N132=-11030.4284
  $ws.Range("H132").Value = 3644.2917
  $ws.Range("J132").Value = 1990.1428
  $ws.Range("L132").Value = 5970.428400000001
  $ws.Range("N132").Value = -11030.4284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
  # row 22: H22=2838.5, I22=860.3333, J22=3431.95, K22=860.3333, L22=3431.95, M22=-565.3333, N22=-4021.95
  $ws.Range("H22").Value = 2838.5
  $ws.Range("I22").Value = 860.3333
  $ws.Range("J22").Value = 3431.95
  $ws.Range("K22").Value = 860.3333
  $ws.Range("L22").Value = 3431.95
  $ws.Range("M22").Value = -565.3333
  $ws.Range("N22").Value = -4021.95
  # row 27: H27=2838.5, I27=860.3333, J27=3431.95, K27=860.3333, L27=3431.95, M27=-753.3333, N27=-3645.95
  $ws.Range("H27").Value = 2838.5
  $ws.Range("I27").Value = 860.3333
  $ws.Range("J27").Value = 3431.95
  $ws.Range("K27").Value = 860.3333
  $ws.Range("L27").Value = 3431.95
  $ws.Range("M27").Value = -753.3333
  $ws.Range("N27").Value = -3645.95
  # row 46: H46=13167.667, I46=5749.25, K46=5749.25, M46=-5561.25
  $ws.Range("H46").Value = 13167.667
  $ws.Range("I46").Value = 5749.25
  $ws.Range("K46").Value = 5749.25
  $ws.Range("M46").Value = -5561.25
  # row 55: H55=1545.8334, I55=741.7143, K55=741.7143, M55=-568.7143
  $ws.Range("H55").Value = 1545.8334
  $ws.Range("I55").Value = 741.7143
  $ws.Range("K55").Value = 741.7143
  $ws.Range("M55").Value = -568.7143
  # row 101: H101=16361, J101=16361, L101=16361, N101=-22851
  $ws.Range("H101").Value = 16361
  $ws.Range("J101").Value = 16361
  $ws.Range("L101").Value = 16361
  $ws.Range("N101").Value = -22851
  # row 136: H136=2936.1765, I136=2740.776, K136=8222.328, M136=-5672.328
  $ws.Range("H136").Value = 2936.1765
  $ws.Range("I136").Value = 2740.776
  $ws.Range("K136").Value = 8222.328
  $ws.Range("M136").Value = -5672.328

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
  # row 81: H81=4797.7896, I81=1969.3636, J81=8686.875, K81=3938.7272, L81=17373.75, M81=-2877.7272, N81=-19495.75
  $ws.Range("H81").Value = 4797.7896
  $ws.Range("I81").Value = 1969.3636
  $ws.Range("J81").Value = 8686.875
  $ws.Range("K81").Value = 3938.7272
  $ws.Range("L81").Value = 17373.75
  $ws.Range("M81").Value = -2877.7272
  $ws.Range("N81").Value = -19495.75
  # row 84: H84=4797.7896, I84=1969.3636, J84=8686.875, K84=19693.636, L84=86868.75, M84=-14389.636, N84=-97476.75
  $ws.Range("H84").Value = 4797.7896
  $ws.Range("I84").Value = 1969.3636
  $ws.Range("J84").Value = 8686.875
  $ws.Range("K84").Value = 19693.636
  $ws.Range("L84").Value = 86868.75
  $ws.Range("M84").Value = -14389.636
  $ws.Range("N84").Value = -97476.75
  # row 103: H103=100000, J103=100000, L103=100000, N103=-102344
  $ws.Range("H103").Value = 100000
  $ws.Range("J103").Value = 100000
  $ws.Range("L103").Value = 100000
  $ws.Range("N103").Value = -102344
  # row 132: H132=1556.3939, I132=1432.0741, K132=4296.2223, M132=-1766.2223
  $ws.Range("H132").Value = 1556.3939
  $ws.Range("I132").Value = 1432.0741
  $ws.Range("K132").Value = 4296.2223
  $ws.Range("M132").Value = -1766.2223
  # row 136: H136=3010.2173, I136=1556.5294, K136=4669.5882, M136=-2119.5882
  $ws.Range("H136").Value = 3010.2173
  $ws.Range("I136").Value = 1556.5294
  $ws.Range("K136").Value = 4669.5882
  $ws.Range("M136").Value = -2119.5882

Write-Host "Applied all Faerie_Profits updates"